# Finished unit testing for Database class
# Updates the "Database ProcessQuery" sheet's test-case table:
#  - removes the "null query and null params", "empty query with params",
#    "empty query with empty params" and "valid query with empty params" rows
#  - adds a new "empty query without params" row in their place
#  - renumbers the trailing "valid query with invalid params" row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Database ProcessQuery")
$ws.Activate()

# Remove the four obsolete test-case rows (old rows 11-14) and open a single
# fresh row in their place for the replacement test case.
$ws.Rows("11:14").Delete()
$ws.Rows("11:11").Insert()

# New row 11: "empty query without params"
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "invalid"
$ws.Cells.Item(11, 3).Value = "empty query without params"
$ws.Cells.Item(11, 4).Value = """"""
$ws.Cells.Item(11, 5).Value = "NA"
$ws.Cells.Item(11, 6).Value = "PDO Exception"

# Row 12 (previously row 15): renumber the "Test Case" id to keep it sequential
$ws.Cells.Item(12, 1).Value = 10

# Restore the view: scrolled so column C is left-most, with D14 selected
# (mirrors the author's on-screen state when they saved).
$ws.Range("D14").Select()

Write-Output "Database ProcessQuery test cases updated"
